$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 2149.5
$ws.Cells.Item(18, 10).Value = 2999
$ws.Cells.Item(18, 12).Value = 2999
$ws.Cells.Item(18, 14).Value = -3567
$ws.Cells.Item(19, 8).Value = 1188
$ws.Cells.Item(19, 9).Value = 1599.6
$ws.Cells.Item(19, 10).Value = 894
$ws.Cells.Item(19, 11).Value = 1599.6
$ws.Cells.Item(19, 12).Value = 894
$ws.Cells.Item(19, 13).Value = -1424.6
$ws.Cells.Item(19, 14).Value = -1244
$ws.Cells.Item(26, 8).Value = 10000
$ws.Cells.Item(26, 10).Value = 10000
$ws.Cells.Item(26, 12).Value = 10000
$ws.Cells.Item(26, 14).Value = -10688
$ws.Cells.Item(28, 8).Value = 126215.5
$ws.Cells.Item(28, 9).Value = 167745.92
$ws.Cells.Item(28, 11).Value = 167745.92
$ws.Cells.Item(28, 13).Value = -167260.92
$ws.Cells.Item(33, 8).Value = 529.9167
$ws.Cells.Item(33, 9).Value = 350.9091
$ws.Cells.Item(33, 11).Value = 350.9091
$ws.Cells.Item(33, 13).Value = -121.9091
$ws.Cells.Item(38, 8).Value = 1712
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 14).ClearContents()
$ws.Cells.Item(41, 8).Value = 2196.7273
$ws.Cells.Item(41, 10).Value = 4016.6667
$ws.Cells.Item(41, 12).Value = 4016.6667
$ws.Cells.Item(41, 14).Value = -4896.6667
$ws.Cells.Item(53, 8).Value = 1060.5186
$ws.Cells.Item(53, 9).Value = 634.5294
$ws.Cells.Item(53, 10).Value = 1784.7
$ws.Cells.Item(53, 11).Value = 634.5294
$ws.Cells.Item(53, 12).Value = 1784.7
$ws.Cells.Item(53, 13).Value = 2.47059999999999
$ws.Cells.Item(53, 14).Value = -3058.7
$ws.Cells.Item(62, 8).Value = 7991.625
$ws.Cells.Item(62, 9).Value = 7062
$ws.Cells.Item(62, 10).Value = 10036.8
$ws.Cells.Item(62, 11).Value = 7062
$ws.Cells.Item(62, 12).Value = 10036.8
$ws.Cells.Item(62, 13).Value = -6438
$ws.Cells.Item(62, 14).Value = -11284.8
$ws.Cells.Item(64, 8).Value = 7828.143
$ws.Cells.Item(64, 10).Value = 8999.799999999999
$ws.Cells.Item(64, 12).Value = 8999.799999999999
$ws.Cells.Item(64, 14).Value = -9495.799999999999
$ws.Cells.Item(65, 8).Value = 7991.625
$ws.Cells.Item(65, 9).Value = 7062
$ws.Cells.Item(65, 10).Value = 10036.8
$ws.Cells.Item(65, 11).Value = 35310
$ws.Cells.Item(65, 12).Value = 50184
$ws.Cells.Item(65, 13).Value = -32190
$ws.Cells.Item(65, 14).Value = -56424
$ws.Cells.Item(67, 8).Value = 7828.143
$ws.Cells.Item(67, 10).Value = 8999.799999999999
$ws.Cells.Item(67, 12).Value = 8999.799999999999
$ws.Cells.Item(67, 14).Value = -10715.8
$ws.Cells.Item(74, 8).Value = 12598.6
$ws.Cells.Item(74, 9).Value = 14331
$ws.Cells.Item(74, 10).Value = 10000
$ws.Cells.Item(74, 11).Value = 14331
$ws.Cells.Item(74, 12).Value = 10000
$ws.Cells.Item(74, 13).Value = -13395
$ws.Cells.Item(74, 14).Value = -11872
$ws.Cells.Item(76, 8).Value = 4233.3335
$ws.Cells.Item(76, 9).Value = 4250
$ws.Cells.Item(76, 11).Value = 4250
$ws.Cells.Item(76, 13).Value = -3935
$ws.Cells.Item(77, 8).Value = 12598.6
$ws.Cells.Item(77, 9).Value = 14331
$ws.Cells.Item(77, 10).Value = 10000
$ws.Cells.Item(77, 11).Value = 71655
$ws.Cells.Item(77, 12).Value = 50000
$ws.Cells.Item(77, 13).Value = -66975
$ws.Cells.Item(77, 14).Value = -59360
$ws.Cells.Item(79, 8).Value = 4233.3335
$ws.Cells.Item(79, 9).Value = 4250
$ws.Cells.Item(79, 11).Value = 4250
$ws.Cells.Item(79, 13).Value = -3158
$ws.Cells.Item(86, 8).Value = 321575140
$ws.Cells.Item(86, 9).Value = 333335330
$ws.Cells.Item(86, 10).Value = 312755000
$ws.Cells.Item(86, 11).Value = 333335330
$ws.Cells.Item(86, 12).Value = 312755000
$ws.Cells.Item(86, 13).Value = -333334207
$ws.Cells.Item(86, 14).Value = -312757246
$ws.Cells.Item(89, 8).Value = 321575140
$ws.Cells.Item(89, 9).Value = 333335330
$ws.Cells.Item(89, 10).Value = 312755000
$ws.Cells.Item(89, 11).Value = 1666676650
$ws.Cells.Item(89, 12).Value = 1563775000
$ws.Cells.Item(89, 13).Value = -1666671034
$ws.Cells.Item(89, 14).Value = -1563786232
$ws.Cells.Item(92, 8).Value = 41668120
$ws.Cells.Item(92, 9).Value = 41668120
$ws.Cells.Item(92, 11).Value = 41668120
$ws.Cells.Item(92, 13).Value = -41666872
$ws.Cells.Item(94, 8).Value = 1838.2858
$ws.Cells.Item(94, 9).Value = 1838.2858
$ws.Cells.Item(94, 11).Value = 1838.2858
$ws.Cells.Item(94, 13).Value = -1387.2858
$ws.Cells.Item(100, 8).Value = 3859
$ws.Cells.Item(100, 9).Value = 2349.5
$ws.Cells.Item(100, 10).Value = 4194.4443
$ws.Cells.Item(100, 11).Value = 2349.5
$ws.Cells.Item(100, 12).Value = 4194.4443
$ws.Cells.Item(100, 13).Value = -1808.5
$ws.Cells.Item(100, 14).Value = -5276.4443
$ws.Cells.Item(101, 8).Value = 1352.7
$ws.Cells.Item(101, 9).Value = 1319.7142
$ws.Cells.Item(101, 11).Value = 3959.1426
$ws.Cells.Item(101, 13).Value = -2337.1426
$ws.Cells.Item(121, 8).Value = 2693.1177
$ws.Cells.Item(121, 10).Value = 2798.9375
$ws.Cells.Item(121, 12).Value = 8396.8125
$ws.Cells.Item(121, 14).Value = -11890.8125
$ws.Cells.Item(131, 8).Value = 2054.611
$ws.Cells.Item(131, 9).Value = 584.6429000000001
$ws.Cells.Item(131, 11).Value = 1753.9287
$ws.Cells.Item(131, 13).Value = 3286.0713
$ws.Cells.Item(132, 8).Value = 2964.082
$ws.Cells.Item(132, 9).Value = 2980.15
$ws.Cells.Item(132, 11).Value = 8940.450000000001
$ws.Cells.Item(132, 13).Value = -6410.450000000001
$ws.Cells.Item(137, 8).Value = 31697.785
$ws.Cells.Item(137, 9).Value = 41594.43
$ws.Cells.Item(137, 11).Value = 124783.29
$ws.Cells.Item(137, 13).Value = -122233.29
$ws.Cells.Item(138, 8).Value = 3169.51
$ws.Cells.Item(138, 9).Value = 1104.28
$ws.Cells.Item(138, 10).Value = 3857.92
$ws.Cells.Item(138, 11).Value = 3312.84
$ws.Cells.Item(138, 12).Value = 11573.76
$ws.Cells.Item(138, 13).Value = 1827.16
$ws.Cells.Item(138, 14).Value = -21853.76
$ws.Cells.Item(141, 8).Value = 530.5
$ws.Cells.Item(141, 9).Value = 475.33334
$ws.Cells.Item(141, 10).Value = 696
$ws.Cells.Item(141, 11).Value = 1426.00002
$ws.Cells.Item(141, 12).Value = 2088
$ws.Cells.Item(141, 13).Value = 3753.99998
$ws.Cells.Item(141, 14).Value = -12448

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12098739
$ws.Cells.Item(32, 9).Value = 12555098
$ws.Cells.Item(32, 10).Value = 7940806.5
$ws.Cells.Item(32, 11).Value = 12555098
$ws.Cells.Item(32, 12).Value = 7940806.5
$ws.Cells.Item(32, 13).Value = -12554811
$ws.Cells.Item(32, 14).Value = -7941380.5
$ws.Cells.Item(61, 8).Value = 2677.3057
$ws.Cells.Item(61, 9).Value = 2575.5
$ws.Cells.Item(61, 10).Value = 2880.9167
$ws.Cells.Item(61, 11).Value = 2575.5
$ws.Cells.Item(61, 12).Value = 2880.9167
$ws.Cells.Item(61, 13).Value = -2363.5
$ws.Cells.Item(61, 14).Value = -3304.9167
$ws.Cells.Item(74, 8).Value = 2816.7222
$ws.Cells.Item(74, 9).Value = 2227.0667
$ws.Cells.Item(74, 10).Value = 5765
$ws.Cells.Item(74, 11).Value = 2227.0667
$ws.Cells.Item(74, 12).Value = 5765
$ws.Cells.Item(74, 13).Value = -1353.0667
$ws.Cells.Item(74, 14).Value = -7513
$ws.Cells.Item(77, 8).Value = 2816.7222
$ws.Cells.Item(77, 9).Value = 2227.0667
$ws.Cells.Item(77, 10).Value = 5765
$ws.Cells.Item(77, 11).Value = 11135.3335
$ws.Cells.Item(77, 12).Value = 28825
$ws.Cells.Item(77, 13).Value = -6767.333499999999
$ws.Cells.Item(77, 14).Value = -37561
$ws.Cells.Item(92, 8).Value = 68000
$ws.Cells.Item(92, 10).Value = 68000
$ws.Cells.Item(92, 12).Value = 68000
$ws.Cells.Item(92, 14).Value = -72992
$ws.Cells.Item(97, 8).Value = 675.13635
$ws.Cells.Item(97, 9).Value = 523.8421
$ws.Cells.Item(97, 11).Value = 523.8421
$ws.Cells.Item(97, 13).Value = -27.84209999999996
$ws.Cells.Item(109, 8).Value = 60869.25
$ws.Cells.Item(109, 10).Value = 60869.25
$ws.Cells.Item(109, 12).Value = 60869.25
$ws.Cells.Item(109, 14).Value = -63643.25
$ws.Cells.Item(122, 8).Value = 5085.476
$ws.Cells.Item(122, 9).Value = 2936
$ws.Cells.Item(122, 10).Value = 7449.9
$ws.Cells.Item(122, 11).Value = 8808
$ws.Cells.Item(122, 12).Value = 22349.7
$ws.Cells.Item(122, 13).Value = -6358
$ws.Cells.Item(122, 14).Value = -27249.7
$ws.Cells.Item(129, 8).Value = 108798.86
$ws.Cells.Item(129, 10).Value = 108798.86
$ws.Cells.Item(129, 12).Value = 108798.86
$ws.Cells.Item(129, 14).Value = -118798.86
$ws.Cells.Item(132, 8).Value = 4811.533
$ws.Cells.Item(132, 9).Value = 5270.25
$ws.Cells.Item(132, 10).Value = 4287.2856
$ws.Cells.Item(132, 11).Value = 15810.75
$ws.Cells.Item(132, 12).Value = 12861.8568
$ws.Cells.Item(132, 13).Value = -13280.75
$ws.Cells.Item(132, 14).Value = -17921.8568
$ws.Cells.Item(136, 8).Value = 2677.3057
$ws.Cells.Item(136, 9).Value = 2575.5
$ws.Cells.Item(136, 10).Value = 2880.9167
$ws.Cells.Item(136, 11).Value = 7726.5
$ws.Cells.Item(136, 12).Value = 8642.750100000001
$ws.Cells.Item(136, 13).Value = -5176.5
$ws.Cells.Item(136, 14).Value = -13742.7501

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1795.3478
$ws.Cells.Item(20, 9).Value = 1398
$ws.Cells.Item(20, 11).Value = 1398
$ws.Cells.Item(20, 13).Value = -1151
$ws.Cells.Item(56, 8).Value = 23999.889
$ws.Cells.Item(56, 9).Value = 23999
$ws.Cells.Item(56, 11).Value = 23999
$ws.Cells.Item(56, 13).Value = -23260
$ws.Cells.Item(82, 8).Value = 34482.832
$ws.Cells.Item(82, 9).Value = 15579.4
$ws.Cells.Item(82, 10).Value = 129000
$ws.Cells.Item(82, 11).Value = 15579.4
$ws.Cells.Item(82, 12).Value = 129000
$ws.Cells.Item(82, 13).Value = -15196.4
$ws.Cells.Item(82, 14).Value = -129766
$ws.Cells.Item(85, 8).Value = 34482.832
$ws.Cells.Item(85, 9).Value = 15579.4
$ws.Cells.Item(85, 10).Value = 129000
$ws.Cells.Item(85, 11).Value = 15579.4
$ws.Cells.Item(85, 12).Value = 129000
$ws.Cells.Item(85, 13).Value = -14253.4
$ws.Cells.Item(85, 14).Value = -131652
$ws.Cells.Item(86, 8).Value = 2328.25
$ws.Cells.Item(86, 9).Value = 1928.8
$ws.Cells.Item(86, 10).Value = 2994
$ws.Cells.Item(86, 11).Value = 1928.8
$ws.Cells.Item(86, 12).Value = 2994
$ws.Cells.Item(86, 13).Value = -805.8
$ws.Cells.Item(86, 14).Value = -5240
$ws.Cells.Item(89, 8).Value = 2328.25
$ws.Cells.Item(89, 9).Value = 1928.8
$ws.Cells.Item(89, 10).Value = 2994
$ws.Cells.Item(89, 11).Value = 9644
$ws.Cells.Item(89, 12).Value = 14970
$ws.Cells.Item(89, 13).Value = -4028
$ws.Cells.Item(89, 14).Value = -26202
$ws.Cells.Item(99, 8).Value = 2691.074
$ws.Cells.Item(99, 9).Value = 1446.8
$ws.Cells.Item(99, 11).Value = 1446.8
$ws.Cells.Item(99, 13).Value = 51.20000000000005
$ws.Cells.Item(105, 8).Value = 2304.1785
$ws.Cells.Item(105, 9).Value = 2059.0908
$ws.Cells.Item(105, 11).Value = 2059.0908
$ws.Cells.Item(105, 13).Value = -312.0907999999999
$ws.Cells.Item(116, 8).Value = 117000
$ws.Cells.Item(116, 10).Value = 117000
$ws.Cells.Item(116, 12).Value = 117000
$ws.Cells.Item(116, 14).Value = -126178
$ws.Cells.Item(132, 8).Value = 120000
$ws.Cells.Item(132, 10).Value = 120000
$ws.Cells.Item(132, 12).Value = 120000
$ws.Cells.Item(132, 14).Value = -130120
$ws.Cells.Item(134, 8).Value = 2167125
$ws.Cells.Item(134, 9).Value = 3969920
$ws.Cells.Item(134, 10).Value = 3771.0667
$ws.Cells.Item(134, 11).Value = 11909760
$ws.Cells.Item(134, 12).Value = 11313.2001
$ws.Cells.Item(134, 13).Value = -11907225
$ws.Cells.Item(134, 14).Value = -16383.2001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2925
$ws.Cells.Item(16, 9).Value = 1300
$ws.Cells.Item(16, 10).Value = 3466.6667
$ws.Cells.Item(16, 11).Value = 1300
$ws.Cells.Item(16, 12).Value = 3466.6667
$ws.Cells.Item(16, 13).Value = -1013
$ws.Cells.Item(16, 14).Value = -4040.6667
$ws.Cells.Item(31, 8).Value = 2905.239
$ws.Cells.Item(31, 9).Value = 1973.8
$ws.Cells.Item(31, 11).Value = 1973.8
$ws.Cells.Item(31, 13).Value = -1678.8
$ws.Cells.Item(34, 8).Value = 2905.239
$ws.Cells.Item(34, 9).Value = 1973.8
$ws.Cells.Item(34, 11).Value = 1973.8
$ws.Cells.Item(34, 13).Value = -1771.8
$ws.Cells.Item(58, 8).Value = 2629.1875
$ws.Cells.Item(58, 10).Value = 3260.5
$ws.Cells.Item(58, 12).Value = 3260.5
$ws.Cells.Item(58, 14).Value = -3666.5
$ws.Cells.Item(86, 8).Value = 45208.805
$ws.Cells.Item(86, 9).Value = 45501.2
$ws.Cells.Item(86, 10).Value = 45161.645
$ws.Cells.Item(86, 11).Value = 45501.2
$ws.Cells.Item(86, 12).Value = 45161.645
$ws.Cells.Item(86, 13).Value = -44378.2
$ws.Cells.Item(86, 14).Value = -47407.645
$ws.Cells.Item(89, 8).Value = 45208.805
$ws.Cells.Item(89, 9).Value = 45501.2
$ws.Cells.Item(89, 10).Value = 45161.645
$ws.Cells.Item(89, 11).Value = 227506
$ws.Cells.Item(89, 12).Value = 225808.225
$ws.Cells.Item(89, 13).Value = -221890
$ws.Cells.Item(89, 14).Value = -237040.225
$ws.Cells.Item(105, 8).Value = 1742.5
$ws.Cells.Item(105, 9).Value = 1449.5834
$ws.Cells.Item(105, 11).Value = 1449.5834
$ws.Cells.Item(105, 13).Value = 297.4166
$ws.Cells.Item(113, 8).Value = 2925
$ws.Cells.Item(113, 9).Value = 1300
$ws.Cells.Item(113, 10).Value = 3466.6667
$ws.Cells.Item(113, 11).Value = 1300
$ws.Cells.Item(113, 12).Value = 3466.6667
$ws.Cells.Item(113, 13).Value = 870
$ws.Cells.Item(113, 14).Value = -7806.6667
$ws.Cells.Item(122, 8).Value = 14290814
$ws.Cells.Item(122, 9).Value = 16672200
$ws.Cells.Item(122, 10).Value = 2500
$ws.Cells.Item(122, 11).Value = 50016600
$ws.Cells.Item(122, 12).Value = 7500
$ws.Cells.Item(122, 13).Value = -50014150
$ws.Cells.Item(122, 14).Value = -12400
$ws.Cells.Item(132, 8).Value = 5627.5625
$ws.Cells.Item(132, 9).Value = 5760.0713
$ws.Cells.Item(132, 10).Value = 4700
$ws.Cells.Item(132, 11).Value = 17280.2139
$ws.Cells.Item(132, 12).Value = 14100
$ws.Cells.Item(132, 13).Value = -14750.2139
$ws.Cells.Item(132, 14).Value = -19160
$ws.Cells.Item(134, 8).Value = 2278.6765
$ws.Cells.Item(134, 9).Value = 1952.0834
$ws.Cells.Item(134, 11).Value = 5856.2502
$ws.Cells.Item(134, 13).Value = -3321.2502
$ws.Cells.Item(136, 8).Value = 2629.1875
$ws.Cells.Item(136, 10).Value = 3260.5
$ws.Cells.Item(136, 12).Value = 9781.5
$ws.Cells.Item(136, 14).Value = -14881.5
$ws.Cells.Item(137, 8).Value = 30000
$ws.Cells.Item(137, 9).Value = 30000
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 11).Value = 30000
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 13).Value = -24900
$ws.Cells.Item(137, 14).ClearContents()
$ws.Cells.Item(140, 8).Value = 403331.34
$ws.Cells.Item(140, 10).Value = 549997
$ws.Cells.Item(140, 12).Value = 549997
$ws.Cells.Item(140, 14).Value = -560357
$ws.Cells.Item(141, 8).Value = 430849.44
$ws.Cells.Item(141, 10).Value = 430849.44
$ws.Cells.Item(141, 12).Value = 430849.44
$ws.Cells.Item(141, 14).Value = -441209.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 38.6
$ws.Cells.Item(2, 9).Value = 50
$ws.Cells.Item(2, 11).Value = 300
$ws.Cells.Item(2, 13).Value = -187
$ws.Cells.Item(4, 8).Value = 93054740
$ws.Cells.Item(4, 9).Value = 113947220
$ws.Cells.Item(4, 10).Value = 76141784
$ws.Cells.Item(4, 11).Value = 341841660
$ws.Cells.Item(4, 12).Value = 228425352
$ws.Cells.Item(4, 13).Value = -341841548
$ws.Cells.Item(4, 14).Value = -228425576
$ws.Cells.Item(5, 8).Value = 2178.4
$ws.Cells.Item(5, 9).Value = 1355.5714
$ws.Cells.Item(5, 10).Value = 4098.3335
$ws.Cells.Item(5, 11).Value = 4066.7142
$ws.Cells.Item(5, 12).Value = 12295.0005
$ws.Cells.Item(5, 13).Value = -3954.7142
$ws.Cells.Item(5, 14).Value = -12519.0005
$ws.Cells.Item(26, 8).Value = 45
$ws.Cells.Item(26, 9).Value = 17.5
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 52.5
$ws.Cells.Item(26, 12).Value = 300
$ws.Cells.Item(26, 13).Value = 235.5
$ws.Cells.Item(26, 14).Value = -876
$ws.Cells.Item(47, 8).Value = 693.6
$ws.Cells.Item(47, 9).Value = 717
$ws.Cells.Item(47, 10).Value = 600
$ws.Cells.Item(47, 11).Value = 2151
$ws.Cells.Item(47, 12).Value = 1800
$ws.Cells.Item(47, 13).Value = -1720
$ws.Cells.Item(47, 14).Value = -2662
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 13).ClearContents()
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 13).ClearContents()
$ws.Cells.Item(113, 8).Value = 996.3333
$ws.Cells.Item(113, 9).Value = 399.25
$ws.Cells.Item(113, 11).Value = 1197.75
$ws.Cells.Item(113, 13).Value = 972.25
$ws.Cells.Item(114, 8).Value = 2999.5715
$ws.Cells.Item(114, 9).Value = 161
$ws.Cells.Item(114, 11).Value = 483
$ws.Cells.Item(114, 13).Value = 2771
$ws.Cells.Item(122, 8).Value = 697979.3
$ws.Cells.Item(122, 9).Value = 1049.5
$ws.Cells.Item(122, 10).Value = 1743374
$ws.Cells.Item(122, 11).Value = 9445.5
$ws.Cells.Item(122, 12).Value = 15690366
$ws.Cells.Item(122, 13).Value = -6995.5
$ws.Cells.Item(122, 14).Value = -15695266
$ws.Cells.Item(131, 8).Value = 1308.9
$ws.Cells.Item(131, 9).Value = 889.8333
$ws.Cells.Item(131, 10).Value = 1937.5
$ws.Cells.Item(131, 11).Value = 2669.4999
$ws.Cells.Item(131, 12).Value = 5812.5
$ws.Cells.Item(131, 13).Value = 2370.5001
$ws.Cells.Item(131, 14).Value = -15892.5
$ws.Cells.Item(132, 8).Value = 1131.1666
$ws.Cells.Item(132, 9).Value = 1224.5
$ws.Cells.Item(132, 11).Value = 11020.5
$ws.Cells.Item(132, 13).Value = -8490.5
$ws.Cells.Item(135, 8).Value = 2178.4
$ws.Cells.Item(135, 9).Value = 1355.5714
$ws.Cells.Item(135, 10).Value = 4098.3335
$ws.Cells.Item(135, 11).Value = 12200.1426
$ws.Cells.Item(135, 12).Value = 36885.0015
$ws.Cells.Item(135, 13).Value = -9665.142600000001
$ws.Cells.Item(135, 14).Value = -41955.0015
$ws.Cells.Item(139, 8).Value = 2575.389
$ws.Cells.Item(139, 9).Value = 2723.8
$ws.Cells.Item(139, 11).Value = 8171.400000000001
$ws.Cells.Item(139, 13).Value = -3031.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(33, 8).Value = 12000
$ws.Cells.Item(33, 9).Value = 12000
$ws.Cells.Item(33, 11).Value = 12000
$ws.Cells.Item(33, 13).Value = -11748
$ws.Cells.Item(52, 8).Value = 38999.668
$ws.Cells.Item(52, 9).Value = 26999
$ws.Cells.Item(52, 11).Value = 26999
$ws.Cells.Item(52, 13).Value = -26740
$ws.Cells.Item(70, 8).Value = 6083.032
$ws.Cells.Item(70, 9).Value = 7734
$ws.Cells.Item(70, 10).Value = 5175
$ws.Cells.Item(70, 11).Value = 7734
$ws.Cells.Item(70, 12).Value = 5175
$ws.Cells.Item(70, 13).Value = -7464
$ws.Cells.Item(70, 14).Value = -5715
$ws.Cells.Item(73, 8).Value = 6083.032
$ws.Cells.Item(73, 9).Value = 7734
$ws.Cells.Item(73, 10).Value = 5175
$ws.Cells.Item(73, 11).Value = 7734
$ws.Cells.Item(73, 12).Value = 5175
$ws.Cells.Item(73, 13).Value = -6798
$ws.Cells.Item(73, 14).Value = -7047
$ws.Cells.Item(80, 8).Value = 2329.8
$ws.Cells.Item(80, 9).Value = 1100
$ws.Cells.Item(80, 10).Value = 3149.6667
$ws.Cells.Item(80, 11).Value = 1100
$ws.Cells.Item(80, 12).Value = 3149.6667
$ws.Cells.Item(80, 13).Value = -102
$ws.Cells.Item(80, 14).Value = -5145.6667
$ws.Cells.Item(83, 8).Value = 2329.8
$ws.Cells.Item(83, 9).Value = 1100
$ws.Cells.Item(83, 10).Value = 3149.6667
$ws.Cells.Item(83, 11).Value = 5500
$ws.Cells.Item(83, 12).Value = 15748.3335
$ws.Cells.Item(83, 13).Value = -508
$ws.Cells.Item(83, 14).Value = -25732.3335
$ws.Cells.Item(97, 8).Value = 987.2105
$ws.Cells.Item(97, 9).Value = 879.17645
$ws.Cells.Item(97, 11).Value = 879.17645
$ws.Cells.Item(97, 13).Value = -383.17645
$ws.Cells.Item(120, 8).Value = 81203.25
$ws.Cells.Item(120, 10).Value = 81203.25
$ws.Cells.Item(120, 12).Value = 81203.25
$ws.Cells.Item(120, 14).Value = -90879.25
$ws.Cells.Item(122, 8).Value = 3750.5454
$ws.Cells.Item(122, 9).Value = 4256.857
$ws.Cells.Item(122, 11).Value = 12770.571
$ws.Cells.Item(122, 13).Value = -10320.571
$ws.Cells.Item(126, 8).Value = 3500
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 10).Value = 3500
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 12).Value = 10500
$ws.Cells.Item(126, 13).ClearContents()
$ws.Cells.Item(126, 14).Value = -15440
$ws.Cells.Item(132, 8).Value = 6033.9375
$ws.Cells.Item(132, 9).Value = 6686.8184
$ws.Cells.Item(132, 11).Value = 20060.4552
$ws.Cells.Item(132, 13).Value = -17530.4552
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 10).Value = 69341.86
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 208025.58
$ws.Cells.Item(136, 13).ClearContents()
$ws.Cells.Item(136, 14).Value = -213125.58

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 862.6667
$ws.Cells.Item(16, 9).Value = 624.5714
$ws.Cells.Item(16, 11).Value = 624.5714
$ws.Cells.Item(16, 13).Value = -454.5714
$ws.Cells.Item(40, 8).Value = 66678268
$ws.Cells.Item(40, 9).Value = 166674180
$ws.Cells.Item(40, 11).Value = 166674180
$ws.Cells.Item(40, 13).Value = -166674044
$ws.Cells.Item(46, 8).Value = 6535.3335
$ws.Cells.Item(46, 9).Value = 7585.4443
$ws.Cells.Item(46, 10).Value = 4960.1665
$ws.Cells.Item(46, 11).Value = 7585.4443
$ws.Cells.Item(46, 12).Value = 4960.1665
$ws.Cells.Item(46, 13).Value = -7397.4443
$ws.Cells.Item(46, 14).Value = -5336.1665
$ws.Cells.Item(68, 8).Value = 3793.0952
$ws.Cells.Item(68, 9).Value = 3115.3635
$ws.Cells.Item(68, 10).Value = 6278.1113
$ws.Cells.Item(68, 11).Value = 3115.3635
$ws.Cells.Item(68, 12).Value = 6278.1113
$ws.Cells.Item(68, 13).Value = -2366.3635
$ws.Cells.Item(68, 14).Value = -7776.1113
$ws.Cells.Item(71, 8).Value = 3793.0952
$ws.Cells.Item(71, 9).Value = 3115.3635
$ws.Cells.Item(71, 10).Value = 6278.1113
$ws.Cells.Item(71, 11).Value = 15576.8175
$ws.Cells.Item(71, 12).Value = 31390.5565
$ws.Cells.Item(71, 13).Value = -11832.8175
$ws.Cells.Item(71, 14).Value = -38878.5565
$ws.Cells.Item(75, 8).Value = 99989
$ws.Cells.Item(75, 10).Value = 99989
$ws.Cells.Item(75, 12).Value = 99989
$ws.Cells.Item(75, 14).Value = -101861
$ws.Cells.Item(78, 8).Value = 99989
$ws.Cells.Item(78, 10).Value = 99989
$ws.Cells.Item(78, 12).Value = 299967
$ws.Cells.Item(78, 14).Value = -309327
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).ClearContents()
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).ClearContents()
$ws.Cells.Item(93, 8).Value = 37038252
$ws.Cells.Item(93, 9).Value = 58824180
$ws.Cells.Item(93, 11).Value = 58824180
$ws.Cells.Item(93, 13).Value = -58822932
$ws.Cells.Item(98, 8).Value = 46663.332
$ws.Cells.Item(98, 10).Value = 46663.332
$ws.Cells.Item(98, 12).Value = 46663.332
$ws.Cells.Item(98, 14).Value = -52653.332
$ws.Cells.Item(100, 8).Value = 3420.375
$ws.Cells.Item(100, 9).Value = 2714.3333
$ws.Cells.Item(100, 10).Value = 4328.143
$ws.Cells.Item(100, 11).Value = 2714.3333
$ws.Cells.Item(100, 12).Value = 4328.143
$ws.Cells.Item(100, 13).Value = -2173.3333
$ws.Cells.Item(100, 14).Value = -5410.143
$ws.Cells.Item(122, 8).Value = 8746.870999999999
$ws.Cells.Item(122, 9).Value = 9014.362999999999
$ws.Cells.Item(122, 11).Value = 27043.089
$ws.Cells.Item(122, 13).Value = -24593.089
$ws.Cells.Item(132, 8).Value = 40296.668
$ws.Cells.Item(132, 9).Value = 48385.207
$ws.Cells.Item(132, 11).Value = 145155.621
$ws.Cells.Item(132, 13).Value = -142625.621
$ws.Cells.Item(136, 8).Value = 2506.9412
$ws.Cells.Item(136, 9).Value = 1742.5714
$ws.Cells.Item(136, 10).Value = 3042
$ws.Cells.Item(136, 11).Value = 5227.7142
$ws.Cells.Item(136, 12).Value = 9126
$ws.Cells.Item(136, 13).Value = -2677.7142
$ws.Cells.Item(136, 14).Value = -14226
$ws.Cells.Item(137, 8).Value = 100000
$ws.Cells.Item(137, 10).Value = 100000
$ws.Cells.Item(137, 12).Value = 100000
$ws.Cells.Item(137, 14).Value = -110200

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 833471.7
$ws.Cells.Item(4, 9).Value = 166
$ws.Cells.Item(4, 10).Value = 5000000
$ws.Cells.Item(4, 11).Value = 166
$ws.Cells.Item(4, 12).Value = 5000000
$ws.Cells.Item(4, 13).Value = -53
$ws.Cells.Item(4, 14).Value = -5000226
$ws.Cells.Item(25, 8).Value = 30027
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 10).Value = 30027
$ws.Cells.Item(25, 11).Value = 0
$ws.Cells.Item(25, 12).Value = 30027
$ws.Cells.Item(25, 13).ClearContents()
$ws.Cells.Item(25, 14).Value = -30613
$ws.Cells.Item(54, 8).Value = 56014
$ws.Cells.Item(54, 10).Value = 62500
$ws.Cells.Item(54, 12).Value = 62500
$ws.Cells.Item(54, 14).Value = -63540
$ws.Cells.Item(96, 8).Value = 12664.429
$ws.Cells.Item(96, 10).Value = 29377
$ws.Cells.Item(96, 12).Value = 29377
$ws.Cells.Item(96, 14).Value = -32123
$ws.Cells.Item(100, 8).Value = 917.65216
$ws.Cells.Item(100, 9).Value = 1100.7059
$ws.Cells.Item(100, 10).Value = 399
$ws.Cells.Item(100, 11).Value = 2201.4118
$ws.Cells.Item(100, 12).Value = 798
$ws.Cells.Item(100, 13).Value = -1660.4118
$ws.Cells.Item(100, 14).Value = -1880
$ws.Cells.Item(107, 8).Value = 463.1875
$ws.Cells.Item(107, 9).Value = 386.5
$ws.Cells.Item(107, 11).Value = 1159.5
$ws.Cells.Item(107, 13).Value = 760.5
$ws.Cells.Item(122, 8).Value = 90913336
$ws.Cells.Item(122, 9).Value = 200000940
$ws.Cells.Item(122, 10).Value = 7008.1665
$ws.Cells.Item(122, 11).Value = 600002820
$ws.Cells.Item(122, 12).Value = 21024.4995
$ws.Cells.Item(122, 13).Value = -600000370
$ws.Cells.Item(122, 14).Value = -25924.4995
$ws.Cells.Item(126, 8).Value = 7093.778
$ws.Cells.Item(126, 9).Value = 8136.5713
$ws.Cells.Item(126, 10).Value = 3444
$ws.Cells.Item(126, 11).Value = 24409.7139
$ws.Cells.Item(126, 12).Value = 10332
$ws.Cells.Item(126, 13).Value = -21939.7139
$ws.Cells.Item(126, 14).Value = -15272
$ws.Cells.Item(132, 8).Value = 2782.544
$ws.Cells.Item(132, 9).Value = 2241.3809
$ws.Cells.Item(132, 10).Value = 4297.8
$ws.Cells.Item(132, 11).Value = 6724.1427
$ws.Cells.Item(132, 12).Value = 12893.4
$ws.Cells.Item(132, 13).Value = -4194.1427
$ws.Cells.Item(132, 14).Value = -17953.4
$ws.Cells.Item(136, 8).Value = 33224.656
$ws.Cells.Item(136, 9).Value = 1459.4286
$ws.Cells.Item(136, 11).Value = 4378.2858
$ws.Cells.Item(136, 13).Value = -1828.2858

Write-Output "Applied 629 cell changes across 8 sheets"